$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The style-number text in column A (rows 43:64) changed from " 8-67G01" to " 867G01"
$ws.Range("A43:A64").Value = " 867G01"

# Update the active selection on the sheet to A43:A64 (active cell A43)
$ws.Range("A43:A64").Select()
